$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "12.664"
$ws.Range("R2").Value = "0.955"
$ws.Range("N3").Value = "0.439"
$ws.Range("N4").Value = "2.59"
$ws.Range("R4").Value = "0.988"
$ws.Range("N5").Value = "4.968"
$ws.Range("R5").Value = "0.245"
$ws.Range("N6").Value = "600.767"
$ws.Range("R6").Value = "0.299"
$ws.Range("N7").Value = "26.695"
$ws.Range("R7").Value = "0.716"
$ws.Range("N8").Value = "31.803"
$ws.Range("R8").Value = "0.812"
$ws.Range("N9").Value = "13.069"
$ws.Range("R9").Value = "0.331"
$ws.Range("N10").Value = "1926.614"
$ws.Range("R10").Value = "0.978"
$ws.Range("N11").Value = "1.979"
$ws.Range("N12").Value = "5.203"
$ws.Range("R12").Value = "0.942"
$ws.Range("N13").Value = "2.021"
$ws.Range("N14").Value = "2.094"
$ws.Range("R14").Value = "0.985"
$ws.Range("N15").Value = "7.093"
$ws.Range("N16").Value = "2.187"
$ws.Range("R16").Value = "0.995"
$ws.Range("N17").Value = "5.264"
$ws.Range("R17").Value = "0.974"
$ws.Range("N18").Value = "5.438"
$ws.Range("N19").Value = "5.025"
$ws.Range("R19").Value = "0.971"
$ws.Range("N20").Value = "3.895"
$ws.Range("R20").Value = "0.868"
$ws.Range("N21").Value = "4.441"
$ws.Range("R21").Value = "0.888"
$ws.Range("N22").Value = "2.595"
$ws.Range("N23").Value = "5.487"
$ws.Range("N24").Value = "4.825"
$ws.Range("N25").Value = "4.773"
$ws.Range("N26").Value = "4.375"
$ws.Range("R26").Value = "0.985"
$ws.Range("N27").Value = "5.162"
$ws.Range("N28").Value = "4.762"
$ws.Range("R28").Value = "0.982"
$ws.Range("N29").Value = "5.003"
$ws.Range("R29").Value = "0.964"
$ws.Range("N30").Value = "4.274"
$ws.Range("R30").Value = "0.979"
$ws.Range("N31").Value = "5.449"
$ws.Range("N32").Value = "6.932"
$ws.Range("R32").Value = "0.991"
$ws.Range("N33").Value = "4.283"
$ws.Range("N34").Value = "5.162"
$ws.Range("N35").Value = "6.21"
$ws.Range("R35").Value = "0.992"
$ws.Range("N36").Value = "3.436"
$ws.Range("N37").Value = "4.433"
$ws.Range("R37").Value = "0.982"
$ws.Range("N38").Value = "3.41"
$ws.Range("N39").Value = "4.681"
$ws.Range("N40").Value = "5.405"
$ws.Range("R40").Value = "0.872"
$ws.Range("N41").Value = "0.269"
$ws.Range("R41").Value = "0.972"
$ws.Range("N42").Value = "0.907"
$ws.Range("R42").Value = "0.995"
$ws.Range("N43").Value = "0.312"
$ws.Range("R43").Value = "0.645"
$ws.Range("N44").Value = "0.448"
$ws.Range("R44").Value = "0.975"
$ws.Range("N45").Value = "0.106"
$ws.Range("R45").Value = "0.992"
$ws.Range("N47").Value = "0.947"
$ws.Range("N51").Value = "0.012"
$ws.Range("N52").Value = "0.019"
$ws.Range("N53").Value = "0.095"
$ws.Range("N54").Value = "0.014"
$ws.Range("N56").Value = "0.903"
$ws.Range("N57").Value = "0.022"
$ws.Range("N58").Value = "0.019"
$ws.Range("N59").Value = "0.014"
$ws.Range("N60").Value = "0.001"
$ws.Range("N61").Value = "0.019"
$ws.Range("N62").Value = "0.011"
$ws.Range("N63").Value = "0.15"
$ws.Range("R63").Value = "0.977"
$ws.Range("N64").Value = "0.757"
$ws.Range("R64").Value = "0.877"
$ws.Range("N65").Value = "0.86"
$ws.Range("N66").Value = "0.466"
$ws.Range("R66").Value = "0.989"
$ws.Range("N67").Value = "0.294"
$ws.Range("R67").Value = "0.854"
$ws.Range("N68").Value = "0.865"
$ws.Range("R68").Value = "0.991"
$ws.Range("N69").Value = "0.315"
$ws.Range("N70").Value = "0.611"
$ws.Range("R70").Value = "0.992"
$ws.Range("N71").Value = "0.086"
$ws.Range("R71").Value = "0.383"
$ws.Range("N72").Value = "0.538"
$ws.Range("R72").Value = "0.981"
$ws.Range("N73").Value = "0.127"
$ws.Range("R73").Value = "0.271"
$ws.Range("N74").Value = "0.273"
$ws.Range("R74").Value = "0.971"
$ws.Range("N75").Value = "0.111"
$ws.Range("R75").Value = "0.27"
$ws.Range("R76").Value = "0.988"
$ws.Range("N77").Value = "0.057"
$ws.Range("R77").Value = "0.878"
